$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the width of column B before inserting, so the new column can match it
$officeColWidth = $ws.Range("B1").EntireColumn.ColumnWidth()

# Insert a new column at C (shifts CLIENT ID..STATUS right by one)
$ws.Range("C1").EntireColumn.Insert()

# New column C inherits the same width as column B
$ws.Range("C1").EntireColumn.ColumnWidth = $officeColWidth

# Add the new header text for the inserted "LOAN OFFICER" column
$ws.Range("C2").Value = "LOAN OFFICER"

# Update selection / view state
$ws.Range("D14").Select() | Out-Null
